# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the crypto
# symbol list with the latest scraped quotes from the GitHub Actions run.
# Values are written as literal text (leading apostrophe forces text entry)
# so the cells keep storing strings like "303.81" / "2.01%" exactly as
# scraped, instead of being reinterpreted as numbers/percentages by Excel.
# Style is reset to "Normal" afterwards so the quote-prefix text format
# that Excel applies on forced-text entry doesn't stick to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'303.81"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'2.01%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'31.88"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'0.61%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'5.186"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'1.57%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.07837"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'4.11%"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'2.366"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'37.92%"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'7.979"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'3.869"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'1.97%"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.9125"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'-1.91%"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.1736"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'2.14%"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.07355"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'-2.36%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.08216"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'3.24%"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'0.03041"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'0.41%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'0.09954"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'0.62%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.001520"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'1.99%"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'0.006003"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'-4.50%"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'3.500"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'1.22%"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'0.88%"
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'0.3244"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'-0.90%"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.1346"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'1.51%"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'4.680"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'2.73%"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'0.04645"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'-0.26%"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'0.29%"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'0.001260"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'3.29%"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'0.004535"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'2.56%"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'3.76%"
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'0.0002740"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'47.48%"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.01795"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'7.45%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.04584"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'1.52%"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.007309"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'3.52%"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'0.1364"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'2.82%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.002239"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'8.65%"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.01100"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'-12.14%"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.00006478"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'7.77%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'0.00000000750"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'-0.07%"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'15.31%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'0.009892"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'-23.70%"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'0.00002099"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'-0.07%"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'0.0001999"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'0.00%"
$c.Style = "Normal"

